$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as text even when they look numeric
# (e.g. "0.998", "602.36"). Assigning such a string straight to .Value would make
# Excel auto-convert it into a real number, changing the cells stored type.
# To keep these cells as text (matching the source data), force a text number
# format on each such cell first, assign the value, and afterwards restore the
# cell style back to the default "Normal" style so no formatting residue remains.
$numericLookingPriceCells = @("D4","D5","D6","D13","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D37","D38","D39","D40","D42","D43","D44","D45","D47","D48","D49","D50")
foreach ($cellRef in $numericLookingPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.335.56"
$ws.Range("E2").Value = "  -0.67%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.830.15"
$ws.Range("E3").Value = "  +3.64%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.29%  "

# Row 5 - BNB
$ws.Range("D5").Value = "602.36"
$ws.Range("E5").Value = "  -2.10%  "

# Row 6 - Solana
$ws.Range("D6").Value = "172.70"
$ws.Range("E6").Value = "  -2.90%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.835.58"
$ws.Range("E7").Value = "  +3.83%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.78%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.30%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +1.73%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.33%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "39.29"
$ws.Range("E13").Value = "  -0.96%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -0.59%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.469.77"
$ws.Range("E15").Value = "  +3.60%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.839.14"
$ws.Range("E16").Value = "  +3.92%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "69.279.08"
$ws.Range("E17").Value = "  -0.81%  "

# Row 18 - Polkadot
$ws.Range("D18").Value = "7.43"
$ws.Range("E18").Value = "  -1.12%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  -3.57%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "16.39"
$ws.Range("E20").Value = "  +0.43%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "500.74"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "9.57"
$ws.Range("E22").Value = "  +4.83%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.745"
$ws.Range("E23").Value = "  +4.85%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "86.96"
$ws.Range("E24").Value = "  +0.97%  "

# Row 25 - Fetch.AI
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -4.43%  "

# Row 26 - PEPE -> InternetComputer(DFINITY)
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "12.55"
$ws.Range("E26").Value = "  -2.85%  "

# Row 27 - InternetComputer(DFINITY) -> PEPE
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "0.0000139"
$ws.Range("E27").Value = "  +8.29%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "10.28"
$ws.Range("E28").Value = "  -9.65%  "

# Row 29 - Dai
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.20%  "

# Row 30 - ImmutableX
$ws.Range("D30").Value = "2.52"
$ws.Range("E30").Value = "  +3.70%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "2.96"
$ws.Range("E31").Value = "  +2.75%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "32.98"
$ws.Range("E32").Value = "  +9.69%  "

# Row 33 - NEARProtocol
$ws.Range("D33").Value = "7.89"
$ws.Range("E33").Value = "  +0.24%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -0.21%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.50%  "

# Row 36 - Mantle
$ws.Range("E36").Value = "  -1.68%  "

# Row 37 - Filecoin
$ws.Range("D37").Value = "6.02"
$ws.Range("E37").Value = "  +0.08%  "

# Row 38 - Kaspa
$ws.Range("D38").Value = "0.139"
$ws.Range("E38").Value = "  +1.96%  "

# Row 39 - Bittensor
$ws.Range("D39").Value = "463.49"
$ws.Range("E39").Value = "  +8.27%  "

# Row 40 - TheGraph
$ws.Range("D40").Value = "0.330"
$ws.Range("E40").Value = "  -1.67%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -0.11%  "

# Row 42 - OKB
$ws.Range("D42").Value = "49.48"
$ws.Range("E42").Value = "  -0.96%  "

# Row 43 - dogwifhat
$ws.Range("D43").Value = "2.87"
$ws.Range("E43").Value = "  -0.39%  "

# Row 44 - Cosmos
$ws.Range("D44").Value = "8.50"
$ws.Range("E44").Value = "  -0.26%  "

# Row 45 - Arweave
$ws.Range("D45").Value = "42.36"
$ws.Range("E45").Value = "  -5.17%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.891.91"
$ws.Range("E46").Value = "  -1.72%  "

# Row 47 - VeChain
$ws.Range("D47").Value = "0.0359"
$ws.Range("E47").Value = "  +0.26%  "

# Row 48 - Monero -> USDe
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.01%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "27.35"
$ws.Range("E49").Value = "  +0.62%  "

# Row 50 - USDe -> Monero
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "138.80"
$ws.Range("E50").Value = "  +2.06%  "

# Row 51 - ThetaToken
$ws.Range("E51").Value = "  -1.64%  "

# Restore default style on the price cells we temporarily reformatted, so only
# the cell type/value differ from the original file (no lingering number format).
foreach ($cellRef in $numericLookingPriceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
